# Apply "finish updates for RQ1 findings" changes:
#  - Merge the separate H4:K11 breakdown into the RQ1 (B/C) table and remove
#    the now-redundant H:K mini table entirely.
#  - Update RQ1 cause-of-flakiness labels/counts to their final values.
#  - Drop the trailing RQ1 rows that no longer apply (Comparisons / duplicate
#    Algorithmic Flakiness / Incorrect Logic rows 12-14).
#  - Move the active selection to D3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update RQ1: Cause of Flakiness? table (columns B:C) ---------------
# Row 6: "Order of Events" (30) -> "Incorrect Program Logic" (45)
$ws.Range("B6").Value = "Incorrect Program Logic"
$ws.Range("C6").Value = 45

# Row 7: "Concurrency" count 6 -> 16
$ws.Range("C7").Value = 16

# Row 8: "Async Wait" count 21 -> 34
$ws.Range("C8").Value = 34

# Row 9: "Delay" (6) -> "Algorithmic Flakiness" (13)
$ws.Range("B9").Value = "Algorithmic Flakiness"
$ws.Range("C9").Value = 13

# Rows 12-14 no longer hold RQ1 data (Comparisons / Algorithmic Flakiness /
# Incorrect Logic rows were removed).
$ws.Range("B12:C14").ClearContents()

# --- Remove the orphaned H4:K11 mini breakdown table --------------------
$ws.Range("H4:K11").ClearContents()

# --- Update the selected cell/view --------------------------------------
$ws.Range("D3").Select()
